# Rename the "Thermal Neutron Porosity" column header to "Neutron Porosity"
# on the Constants worksheet (cell D1).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")

$ws.Range("D1").Value = "Neutron Porosity"
